# Medical Device Report.pptx - apply authored changes
#  1) Date placeholder text on the slide master + every slide layout:
#     10/27/2020 -> 10/28/2020
#  2) Slide 7 ("Medical Device Events In Radiology"):
#       - retitle to "Recalls In Radiology"
#       - remove the added commentary textbox ("Radiology had the highest ...")
#  3) Slide 8 ("Top manufacturers of devices across the world"):
#       - clear the title text (now empty)
#       - resize/reposition the map picture and turn off its fill

$p = $ppt.ActivePresentation

# --- 1) Refresh the cached "10/27/2020" date text everywhere it appears ---
$newDate = "10/28/2020"

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "*Date*" -and $sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "10/27/2020") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "*Date*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "10/27/2020") {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2) Slide 7: "Medical Device Events In Radiology" -> "Recalls In Radiology" ---
$s7 = $p.Slides.Item(7)
for ($j = 1; $j -le $s7.Shapes.Count; $j++) {
    $sh = $s7.Shapes.Item($j)
    if ($sh.Name -eq "Title 1") {
        $sh.TextFrame.TextRange.Text = "Recalls In Radiology"
    }
}

# Remove the "Radiology had the highest amount..." textbox (TextBox 2)
for ($j = $s7.Shapes.Count; $j -ge 1; $j--) {
    $sh = $s7.Shapes.Item($j)
    if ($sh.Name -eq "TextBox 2") {
        $sh.Delete()
    }
}

# --- 3) Slide 8: clear title text; reposition/resize picture; remove its fill ---
$s8 = $p.Slides.Item(8)
for ($j = 1; $j -le $s8.Shapes.Count; $j++) {
    $sh = $s8.Shapes.Item($j)
    if ($sh.Name -eq "Title 1") {
        $sh.TextFrame.TextRange.Text = ""
    } elseif (-not $sh.HasTextFrame) {
        # the map picture placeholder
        $sh.LockAspectRatio = 0
        $sh.Left = 93.85992125984252
        $sh.Top = 165.6
        $sh.Width = 772.2800787401575
        $sh.Height = 303.120047
        $sh.Fill.Visible = 0
    }
}
